$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 556.903067169749
$ws.Range("D2").Value = 137.0797485771566
$ws.Range("F2").Value = 452
$ws.Range("G2").Value = 515
$ws.Range("H2").Value = 629
$ws.Range("C3").Value = 39.7527068601508
$ws.Range("D3").Value = 5.498549166859621
$ws.Range("E3").Value = 19.17
$ws.Range("F3").Value = 36.56
$ws.Range("G3").Value = 39.55
$ws.Range("H3").Value = 42.89
$ws.Range("C4").Value = 1.699852359697954
$ws.Range("D4").Value = 2.200329288054807
$ws.Range("F4").Value = 0.57
$ws.Range("G4").Value = 1.11
$ws.Range("H4").Value = 2.14
$ws.Range("C5").Value = 323.5643501279831
$ws.Range("D5").Value = 11.37798183850485
$ws.Range("F5").Value = 316.09
$ws.Range("G5").Value = 325.18
$ws.Range("H5").Value = 332.46
$ws.Range("C6").Value = 20.86686181265001
$ws.Range("D6").Value = 2.527131286839315
$ws.Range("F6").Value = 19.48
$ws.Range("G6").Value = 20.81
$ws.Range("H6").Value = 22.31
$ws.Range("I6").Value = 37.8
$ws.Range("C7").Value = -76.24924404900628
$ws.Range("D7").Value = 23.01777259675368
$ws.Range("F7").Value = -93
$ws.Range("C8").Value = 7.557583971088439
$ws.Range("D8").Value = 6.967010693220619
$ws.Range("F8").Value = 7.8
$ws.Range("C9").Value = 9.322028731546192
$ws.Range("D9").Value = 1.685296913919198
$ws.Range("C10").Value = 867.8301435397293
$ws.Range("D10").Value = 0.461476958793068
$ws.Range("C11").Value = 0.555614872938282
$ws.Range("D11").Value = 0.5888382736154281
$ws.Range("C12").Value = 22.73680810202462
$ws.Range("D12").Value = 12.29101085861106
$ws.Range("C13").Value = 0.6739719803998502
$ws.Range("D13").Value = 0.7505078519842362
$ws.Range("C14").Value = 1.826418901577393
$ws.Range("D14").Value = 1.663703329637374
$ws.Range("C15").Value = 93.64924404900609
$ws.Range("D15").Value = 23.01777259675369
$ws.Range("H15").Value = 110.4
$ws.Range("C16").Value = -85.44001163749608
$ws.Range("D16").Value = 20.67387923932912
$ws.Range("F16").Value = -102.265723755961
$ws.Range("G16").Value = -83.5175485570292
$ws.Range("H16").Value = -68.7376019773414
$ws.Range("C17").Value = -77.88242766640762
$ws.Range("D17").Value = 25.46912265037626
$ws.Range("F17").Value = -93.39612087980606
$ws.Range("G17").Value = -72.61209675612977
$ws.Range("H17").Value = -57.69982180459142
